$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.18599966666667
$ws.Range("H2").Value = 63.557999
$ws.Range("I2").Value = 0.08765141600314529
$ws.Range("J2").Value = 0.08765141600314529
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 55.04475107194456
$ws.Range("R2").Value = 495.402759647501
$ws.Range("S2").Value = 0.03038009755198353
$ws.Range("T2").Value = 0.03038009755198353
$ws.Range("G3").Value = 21.18599966666667
$ws.Range("H3").Value = 63.557999
$ws.Range("I3").Value = 0.08765141600314529
$ws.Range("J3").Value = 0.08765141600314529
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("N3").Value = 13.00021
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 91.80748157553222
$ws.Range("R3").Value = 826.2673341797899
$ws.Range("S3").Value = 0.05067004922269819
$ws.Range("T3").Value = 0.05067004922269819
$ws.Range("G4").Value = 21.18599966666667
$ws.Range("H4").Value = 63.557999
$ws.Range("I4").Value = 0.08765141600314529
$ws.Range("J4").Value = 0.08765141600314529
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 9.940562849598557
$ws.Range("R4").Value = 89.465065646387
$ws.Range("S4").Value = 0.00548635906623892
$ws.Range("T4").Value = 0.00548635906623892
$ws.Range("G5").Value = 21.18599966666667
$ws.Range("H5").Value = 63.557999
$ws.Range("I5").Value = 0.08765141600314529
$ws.Range("J5").Value = 0.08765141600314529
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 2.020070944216889
$ws.Range("R5").Value = 18.180638497952
$ws.Range("S5").Value = 0.001114910162224639
$ws.Range("T5").Value = 0.001114910162224639
$ws.Range("I6").Value = 0.5040014103551328
$ws.Range("J6").Value = 0.5040014103551328
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 316.5109411570915
$ws.Range("R6").Value = 2848.598470413823
$ws.Range("S6").Value = 0.1746875602372103
$ws.Range("T6").Value = 0.1746875602372103
$ws.Range("I7").Value = 0.5040014103551328
$ws.Range("J7").Value = 0.5040014103551328
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("N7").Value = 13.00021
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("Q7").Value = 527.8990609069078
$ws.Range("R7").Value = 4751.091548162171
$ws.Range("S7").Value = 0.2913561176249279
$ws.Range("T7").Value = 0.2913561176249279
$ws.Range("I8").Value = 0.5040014103551328
$ws.Range("J8").Value = 0.5040014103551328
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 57.15889057333345
$ws.Range("R8").Value = 514.4300151600011
$ws.Range("S8").Value = 0.03154692568799871
$ws.Range("T8").Value = 0.03154692568799871
$ws.Range("I9").Value = 0.5040014103551328
$ws.Range("J9").Value = 0.5040014103551328
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 11.61554087005511
$ws.Range("R9").Value = 104.539867830496
$ws.Range("S9").Value = 0.006410806804995873
$ws.Range("T9").Value = 0.006410806804995872
$ws.Range("G10").Value = 37.20718233333333
$ws.Range("H10").Value = 111.621547
$ws.Range("I10").Value = 0.1539347809079331
$ws.Range("J10").Value = 0.1539347809079331
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 96.67044849666144
$ws.Range("R10").Value = 870.034036469953
$ws.Range("S10").Value = 0.05335400012771507
$ws.Range("T10").Value = 0.05335400012771508
$ws.Range("G11").Value = 37.20718233333333
$ws.Range("H11").Value = 111.621547
$ws.Range("I11").Value = 0.1539347809079331
$ws.Range("J11").Value = 0.1539347809079331
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("N11").Value = 13.00021
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 161.2337279472077
$ws.Range("R11").Value = 1451.10355152487
$ws.Range("S11").Value = 0.08898752902531935
$ws.Range("T11").Value = 0.08898752902531937
$ws.Range("G12").Value = 37.20718233333333
$ws.Range("H12").Value = 111.621547
$ws.Range("I12").Value = 0.1539347809079331
$ws.Range("J12").Value = 0.1539347809079331
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 17.45777118192344
$ws.Range("R12").Value = 157.119940637311
$ws.Range("S12").Value = 0.009635229176599215
$ws.Range("T12").Value = 0.009635229176599217
$ws.Range("G13").Value = 37.20718233333333
$ws.Range("H13").Value = 111.621547
$ws.Range("I13").Value = 0.1539347809079331
$ws.Range("J13").Value = 0.1539347809079331
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 3.547680030695111
$ws.Range("R13").Value = 31.92912027625599
$ws.Range("S13").Value = 0.00195802257829947
$ws.Range("T13").Value = 0.00195802257829947
$ws.Range("G14").Value = 61.49336899999999
$ws.Range("H14").Value = 184.480107
$ws.Range("I14").Value = 0.2544123927337887
$ws.Range("J14").Value = 0.2544123927337887
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 159.7700010590437
$ws.Range("R14").Value = 1437.930009531393
$ws.Range("S14").Value = 0.08817967423833402
$ws.Range("T14").Value = 0.08817967423833402
$ws.Range("G15").Value = 61.49336899999999
$ws.Range("H15").Value = 184.480107
$ws.Range("I15").Value = 0.2544123927337887
$ws.Range("J15").Value = 0.2544123927337887
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("N15").Value = 13.00021
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 266.4755702024966
$ws.Range("R15").Value = 2398.28013182247
$ws.Range("S15").Value = 0.1470722214256403
$ws.Range("T15").Value = 0.1470722214256403
$ws.Range("G16").Value = 61.49336899999999
$ws.Range("H16").Value = 184.480107
$ws.Range("I16").Value = 0.2544123927337887
$ws.Range("J16").Value = 0.2544123927337887
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 28.85295520606567
$ws.Range("R16").Value = 259.676596854591
$ws.Range("S16").Value = 0.01592441743768831
$ws.Range("T16").Value = 0.01592441743768831
$ws.Range("G17").Value = 61.49336899999999
$ws.Range("H17").Value = 184.480107
$ws.Range("I17").Value = 0.2544123927337887
$ws.Range("J17").Value = 0.2544123927337887
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 5.863351738570666
$ws.Range("R17").Value = 52.77016564713599
$ws.Range("S17").Value = 0.003236079632126064
$ws.Range("T17").Value = 0.003236079632126063
